$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Columns A, B and D hold values that look like a date / time / number
# ("2023-06-08", "18:10:32", "23") but must be stored as literal text, just
# like the rest of the table. Prefixing with an apostrophe forces Excel to
# keep them as text instead of auto-converting them to date/time/numeric
# values. Column C ("Thursday") is plain text already.
$ws.Cells.Item($row, 1).Value = "'2023-06-08"
$ws.Cells.Item($row, 2).Value = "'18:10:32"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "'23"

# Re-apply the same formatting used by the rest of the data rows (default,
# General number format, no quote-prefix styling) so the new row's cells
# match the look of the existing ones.
$srcFormat = $ws.Range("A2:D2")
$srcFormat.Copy()
$dstFormat = $ws.Range("A$row`:D$row")
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Numeric columns (Beijing ... Wuhan)
$ws.Cells.Item($row, 5).Value = 119783
$ws.Cells.Item($row, 6).Value = 134315
$ws.Cells.Item($row, 7).Value = 160473
$ws.Cells.Item($row, 8).Value = 131465
$ws.Cells.Item($row, 9).Value = 175683
$ws.Cells.Item($row, 10).Value = 113238
$ws.Cells.Item($row, 11).Value = 201256
$ws.Cells.Item($row, 12).Value = 221248
$ws.Cells.Item($row, 13).Value = 172810
$ws.Cells.Item($row, 14).Value = 120139
$ws.Cells.Item($row, 15).Value = 38620
$ws.Cells.Item($row, 16).Value = 34422
$ws.Cells.Item($row, 17).Value = 50868
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36832
$ws.Cells.Item($row, 20).Value = -1
